{"js": "// Replace the sentence about the added ammeter/conductor/resistor with the\n// sentence about the added voltage source, preserving the run's formatting.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst oldText =\n  \"\u041d\u0430 \u0441\u0445\u0435\u043c\u0443 \u0431\u044b\u043b \u0434\u043e\u0431\u0430\u0432\u043b\u0435\u043d \u0430\u043c\u043f\u0435\u0440\u043c\u0435\u0442\u0440. \u041d\u0430 \u0441\u0445\u0435\u043c\u0443 \u0431\u044b\u043b \u0434\u043e\u0431\u0430\u0432\u043b\u0435\u043d \u043f\u0440\u043e\u0432\u043e\u0434\u043d\u0438\u043a. \u041d\u0430 \u0441\u0445\u0435\u043c\u0443 \u0431\u044b\u043b \u0434\u043e\u0431\u0430\u0432\u043b\u0435\u043d \u0440\u0435\u0437\u0438\u0441\u0442\u043e\u0440. \";\nconst newText = \"\u041d\u0430 \u0441\u0445\u0435\u043c\u0443 \u0431\u044b\u043b \u0434\u043e\u0431\u0430\u0432\u043b\u0435\u043d \u0438\u0441\u0442\u043e\u0447\u043d\u0438\u043a \u043d\u0430\u043f\u0440\u044f\u0436\u0435\u043d\u0438\u044f. \";\n\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  if (para.text && para.text.indexOf(oldText.trim()) !== -1) {\n    target = para;\n    break;\n  }\n}\n\nif (!target) {\n  // Fallback: first non-empty paragraph.\n  for (let i = 0; i < paragraphs.items.length; i++) {\n    if (paragraphs.items[i].text && paragraphs.items[i].text.trim().length > 0) {\n      target = paragraphs.items[i];\n      break;\n    }\n  }\n}\n\n// Replace the paragraph's whole range text in-place so the existing run\n// formatting (font, size, color, language, kerning) is kept.\nconst range = target.getRange();\nrange.insertText(newText, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Replace the sentence about the added ammeter/conductor/resistor with the\n# sentence about the added voltage source, preserving the run's formatting.\n$d = $word.ActiveDocument\n\n$oldText = \"\u041d\u0430 \u0441\u0445\u0435\u043c\u0443 \u0431\u044b\u043b \u0434\u043e\u0431\u0430\u0432\u043b\u0435\u043d \u0430\u043c\u043f\u0435\u0440\u043c\u0435\u0442\u0440. \u041d\u0430 \u0441\u0445\u0435\u043c\u0443 \u0431\u044b\u043b \u0434\u043e\u0431\u0430\u0432\u043b\u0435\u043d \u043f\u0440\u043e\u0432\u043e\u0434\u043d\u0438\u043a. \u041d\u0430 \u0441\u0445\u0435\u043c\u0443 \u0431\u044b\u043b \u0434\u043e\u0431\u0430\u0432\u043b\u0435\u043d \u0440\u0435\u0437\u0438\u0441\u0442\u043e\u0440. \"\n$newText = \"\u041d\u0430 \u0441\u0445\u0435\u043c\u0443 \u0431\u044b\u043b \u0434\u043e\u0431\u0430\u0432\u043b\u0435\u043d \u0438\u0441\u0442\u043e\u0447\u043d\u0438\u043a \u043d\u0430\u043f\u0440\u044f\u0436\u0435\u043d\u0438\u044f. \"\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = $oldText\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = $newText\n$found = $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n\nif (-not $found) {\n    # Fallback: locate the paragraph containing the old sentence and set its\n    # range text directly, which keeps the existing run formatting intact.\n    foreach ($p in $d.Paragraphs) {\n        if ($p.Range.Text -like \"*$oldText*\") {\n            $p.Range.Text = $newText\n            break\n        }\n    }\n}\n"}
